$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The original sheet1 had data occupying columns B:E, with column A left
# completely empty/unnamed. Shift every value one column to the left
# (B->A, C->B, D->C, E->D) to remove that empty leading column.
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 1; $r -le 4; $r++) {
  for ($c = 1; $c -le 4; $c++) {
    $srcCell = $ws.Cells.Item($r, $c + 1)
    $dstCell = $ws.Cells.Item($r, $c)
    $dstCell.Value = $srcCell.Value2
  }
}

# Column E (5) is now vacated (its data moved into D); clear it out.
for ($r = 1; $r -le 4; $r++) {
  $ws.Cells.Item($r, 5).ClearContents()
}

# Remove the now-unused trailing blank sheets
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

Write-Host "Done"
